{"js": "// The sentence about Session Storage persistence needs the phrase\n// \"duration of the browser session\" highlighted (light gray) while the\n// rest of the run stays plain. Find that exact phrase and only change\n// its character formatting \u2014 Word will automatically split the\n// surrounding run(s) around the highlighted portion when the document\n// is saved.\nconst body = context.document.body;\nconst results = body.search(\"duration of the browser session\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the target phrase \"duration of the browser session\".');\n}\n\nconst target = results.items[0];\ntarget.font.highlightColor = \"lightGray\";\nawait context.sync();\n", "ps1": "# The sentence about Session Storage persistence needs the phrase\n# \"duration of the browser session\" highlighted (light gray) while the\n# rest of the sentence keeps its original (default) formatting. Locate\n# that exact phrase with Find and restrict the change to the Font of the\n# found Range only -- Word will split the existing run into three runs\n# (before / highlighted / after) automatically on save.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"duration of the browser session\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute()\nif (-not $found) {\n    throw 'Could not find the target phrase \"duration of the browser session\".'\n}\n\n$rng.Font.HighlightColorIndex = \"wdGray25\"\n"}
